$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Gasera CH4 headers to the CCH4 naming used after the
# transformation (mass-of-carbon basis).
$ws.Range("C1").Value = "avg_Gasera_CCH4_flux_mgm2h"
$ws.Range("D1").Value = "avg_Gasera_CCH4_flux_mgm2h_cor"

# Convert the CH4 flux values (columns C and D, rows 2-18) to the
# CCH4 basis by applying the CH4->C mass-ratio conversion twice
# (12/16 * 12/16 = 0.5625), matching the data transformation applied
# to Master_GHG_2023.
$factor = 0.5625
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    if ($cCell.Value2 -ne $null) {
        $cCell.Value2 = $cCell.Value2 * $factor
    }
    if ($dCell.Value2 -ne $null) {
        $dCell.Value2 = $dCell.Value2 * $factor
    }
}
